# Correction of the name of the variable amount in the database
# (several "Price" values in column E were corrected to 59.39)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (Title) so the longer book titles are fully visible
$ws.Columns.Item(1).ColumnWidth = 66

# Correct the Price (column E) values on the affected rows
$ws.Range("E19").Value = 59.39
$ws.Range("E28").Value = 59.39
$ws.Range("E36").Value = 59.39
$ws.Range("E39").Value = 59.39
$ws.Range("E49").Value = 59.39

# Update the active selection to match the last edited cell
$ws.Range("E19").Select()
